$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.923.87"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.786.08"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "657.23"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.40"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.780.03"
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.458"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("E12").Value = "  +2.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  -4.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.21"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.414.64"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.780.06"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.796.05"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.05"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.74"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.64"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000144"
$ws.Range("E24").Value = "  -4.77%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.28"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.37"
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.930.71"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.72"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.21"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.90"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.176"
$ws.Range("E35").Value = "  +16.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.736.43"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.93"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.32"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.88"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.960"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.20"
$ws.Range("E45").Value = "  +6.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.02"
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.18"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.76"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.39"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.43"
$ws.Range("E51").Value = "  -0.55%  "
